# Final analysis for draft
# Update the "per_day" sheet (sheet2) with the new weekly/trip/km columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("per_day")

# Rename the "day" column header to "week"
$ws.Range("C1").Value = "week"

# Add two new trailing columns: trips (K) and kmd (L)
$ws.Range("K1").Value = "trips"
$ws.Range("L1").Value = "kmd"

# Update the selected range to match the new data extent used for analysis
[void]$ws.Range("A2:L51949").Select()
